# Update lab_info (lympho, neutro, etc...)
#
# Adds 10 new rows (168-177) to the "basic_info_all" sheet describing new
# lab fields (differential WBC counts / percentages), mirroring the
# existing A/B/E column layout:
#   A = "blood" (sort key, reused shared string)
#   B = English label for rows 168-172, internal field code for rows 173-177
#   E = English label (same as B for rows 168-172; descriptive label for 173-177)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that land in column E (cheonan), top to bottom, rows 168-177.
$eValues = @(
    "Lymphocyte",
    "Monocytes",
    "Neutrophil",
    "Eosinophil",
    "Basophil",
    "Atypical Lymphocyte %",
    "Band Neutrophil %",
    "Metamyelocyte %",
    "Myelocyte %",
    "Blast %"
)

# Values that land in column B (vname), top to bottom, rows 168-177.
$bValues = @(
    "Lymphocyte",
    "Monocytes",
    "Neutrophil",
    "Eosinophil",
    "Basophil",
    "aty_lympho_percent",
    "band_neutro_percent",
    "metamyelo_percent",
    "myelo_percent",
    "blast_percent"
)

# Column A ("blood") for every new row - reuses the existing shared string.
for ($i = 0; $i -lt 10; $i++) {
    $row = 168 + $i
    $ws.Range("A$row").Value2 = "blood"
}

# Fill column E first (matches the order new shared strings were authored in).
for ($i = 0; $i -lt 10; $i++) {
    $row = 168 + $i
    $ws.Range("E$row").Value2 = $eValues[$i]
}

# Then column B.
for ($i = 0; $i -lt 10; $i++) {
    $row = 168 + $i
    $ws.Range("B$row").Value2 = $bValues[$i]
}

# Apply the "Consolas 10pt, vertically centered" formatting used by the other
# entries in this block (copy the font from an already-styled cell, then
# drop the horizontal-left override so only vertical centering applies).
$ws.Range("E163").Copy() | Out-Null
for ($i = 0; $i -lt 10; $i++) {
    $row = 168 + $i
    $ws.Range("B$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null
}
for ($i = 0; $i -lt 10; $i++) {
    $row = 168 + $i
    $ws.Range("B$row").HorizontalAlignment = 1
    $ws.Range("E$row").HorizontalAlignment = 1
}

# Update the view/selection to match where the author ended up working.
$ws.Range("C174").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 154
$excel.ActiveWindow.ScrollColumn = 1
